$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old ECs-sending rows (5-7 become obsolete after reshuffle); the MuSCs-sending
# rows (previously 5-7) move up to become rows 2-4 with freshly recalculated TPM-based values.
$ws.Rows("5:7").Delete()

$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Osm"
$ws.Range("C2").Value = "Osmr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03867233333333333
$ws.Range("H2").Value = 0.116017
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 35.22755233333334
$ws.Range("N2").Value = 105.682657
$ws.Range("O2").Value = 0.3133918681326672
$ws.Range("P2").Value = 0.3133918681326672
$ws.Range("Q2").Value = 1.362331646352111
$ws.Range("R2").Value = 12.260984817169
$ws.Range("S2").Value = 0.3133918681326672
$ws.Range("T2").Value = 0.3133918681326672
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Osm"
$ws.Range("C3").Value = "Osmr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03867233333333333
$ws.Range("H3").Value = 0.116017
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 52.87700633333333
$ws.Range("N3").Value = 158.631019
$ws.Range("O3").Value = 0.4704051998635747
$ws.Range("P3").Value = 0.4704051998635747
$ws.Range("Q3").Value = 2.044877214591444
$ws.Range("R3").Value = 18.403894931323
$ws.Range("S3").Value = 0.4704051998635747
$ws.Range("T3").Value = 0.4704051998635747
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Osm"
$ws.Range("C4").Value = "Osmr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03867233333333333
$ws.Range("H4").Value = 0.116017
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.30280066666667
$ws.Range("N4").Value = 72.908402
$ws.Range("O4").Value = 0.2162029320037581
$ws.Range("P4").Value = 0.216202932003758
$ws.Range("Q4").Value = 0.9398460083148887
$ws.Range("R4").Value = 8.458614074833999
$ws.Range("S4").Value = 0.2162029320037581
$ws.Range("T4").Value = 0.216202932003758
